# Generate Report for Handback
# This script fills in the "Latest Target File" (F) and "Latest Handback File" (G)
# columns for the zh-cn and de-de localization status sheets, updates the
# "Status" text from "Ready for handoff" to "Handed back: in sync with en-US"
# everywhere it appears, and records the actual handback timestamps in the
# "Latest Handback DateTime" (H) column.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update the Status text everywhere "Ready for handoff" is used.
# ---------------------------------------------------------------------
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (F) / Latest Handback File (G)
#    with hyperlinked file names, and record the handback datetime (H).
# ---------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8ff0dd60b3ae1acfffdf6aa37c18561b32bc975/e2e/155daaac-08c1-48d0-964f-fe84de3a97db.md",
    [Type]::Missing,
    [Type]::Missing,
    "155daaac-08c1-48d0-964f-fe84de3a97db.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59089f56b5048ba24c9bf369d5af4cdfe414d01e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/155daaac-08c1-48d0-964f-fe84de3a97db.49491f9a83640f5ccaffa0cbee4e071504851fb4.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "155daaac-08c1-48d0-964f-fe84de3a97db.49491f9a83640f5ccaffa0cbee4e071504851fb4.zh-cn.xlf"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8ff0dd60b3ae1acfffdf6aa37c18561b32bc975/e2e/21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.md",
    [Type]::Missing,
    [Type]::Missing,
    "21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59089f56b5048ba24c9bf369d5af4cdfe414d01e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.95050232ed882289749e4ee62910b6e80be082b0.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.95050232ed882289749e4ee62910b6e80be082b0.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("H2").Value = "2016-03-19 16:50:43"
$wsZhCn.Range("H3").Value = "2016-03-19 16:50:43"

# ---------------------------------------------------------------------
# 3. de-de sheet: fill in Latest Target File (F) / Latest Handback File (G)
#    with hyperlinked file names, and record the handback datetime (H).
# ---------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8ff0dd60b3ae1acfffdf6aa37c18561b32bc975/e2e/155daaac-08c1-48d0-964f-fe84de3a97db.md",
    [Type]::Missing,
    [Type]::Missing,
    "155daaac-08c1-48d0-964f-fe84de3a97db.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e4df191ceb0cabbc617e4fc490843ead14c2481/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/155daaac-08c1-48d0-964f-fe84de3a97db.49491f9a83640f5ccaffa0cbee4e071504851fb4.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "155daaac-08c1-48d0-964f-fe84de3a97db.49491f9a83640f5ccaffa0cbee4e071504851fb4.de-de.xlf"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8ff0dd60b3ae1acfffdf6aa37c18561b32bc975/e2e/21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.md",
    [Type]::Missing,
    [Type]::Missing,
    "21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e4df191ceb0cabbc617e4fc490843ead14c2481/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.95050232ed882289749e4ee62910b6e80be082b0.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "21a6145f-f2ad-4ad9-ae1b-10e2d89eec22.95050232ed882289749e4ee62910b6e80be082b0.de-de.xlf"
) | Out-Null

$wsDeDe.Range("H2").Value = "2016-03-19 16:50:49"
$wsDeDe.Range("H3").Value = "2016-03-19 16:50:49"

Write-Output "Handback report generated."
